$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "奶奶烘焙饼干，每个约 +0.6 CPS。"
$ws.Range("I5").Value = 0.6

$ws.Range("F6").Value = "种植饼干树，每个约 +6.5 CPS。"
$ws.Range("I6").Value = 6.5

$ws.Range("F7").Value = "开采饼干矿脉，每个约 +65 CPS。"
$ws.Range("I7").Value = 65

$ws.Range("F8").Value = "批量生产饼干，每个约 +650 CPS。"
$ws.Range("I8").Value = 650

$ws.Range("F9").Value = "资本运作，每个约 +7000 CPS。"
$ws.Range("I9").Value = 7000

$ws.Range("F10").Value = "奶奶祈祷，每个约 +85K CPS。"
$ws.Range("I10").Value = 85000

$ws.Range("F11").Value = "魔法师召唤饼干，每个约 +1.2M CPS。"
$ws.Range("I11").Value = 1200000

$ws.Range("F12").Value = "宇宙运输，每个约 +18M CPS。"
$ws.Range("I12").Value = 18000000

$ws.Range("F13").Value = "跨维度门户，每个约 +260M CPS。"
$ws.Range("I13").Value = 260000000

$ws.Range("F14").Value = "扭曲时间，每个约 +3.2B CPS。"
$ws.Range("I14").Value = 3200000000

$ws.Range("F15").Value = "凝聚物质，每个约 +42B CPS。"
$ws.Range("I15").Value = 42000000000
